$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (shifts existing rows 4-14 down to 5-15)
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row with the new bookmark entry.
# Column A ("ID") holds numeric-looking values that are stored as text
# throughout the sheet, so write it as a text formula and flatten it back
# to a plain value via copy/paste-values (avoids Excel's automatic
# number conversion while keeping the cell's style untouched).
$ws.Cells.Item(4, 1).Formula = "=""5"""
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4163)
$ws.Cells.Item(4, 2).Value = "Java Docs"
$ws.Cells.Item(4, 3).Value = "Oracle.com"
$ws.Cells.Item(4, 4).Value = "Java docs"
$ws.Cells.Item(4, 5).Value = "Coding"
